# Swap the "Lesson(s)" and "Due" entries between the Week 9 row and the
# Week 11 row of the schedule table (Table 1).
#
#   Week 9  (row 10): "Data Wrangling IV" / "Reproducible Report: Initial Analysis"
#   Week 11 (row 12): "Functional Programming" / "Assignment 8"
#
# After the edit these two rows trade their Lesson/Due contents.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Header row is row 1, so Week 9 is row 10 and Week 11 is row 12.
$week9Lesson = $t.Cell(10, 3).Range
$week9Due    = $t.Cell(10, 4).Range
$week11Lesson = $t.Cell(12, 3).Range
$week11Due    = $t.Cell(12, 4).Range

$week9Lesson.Text = "Functional Programming"
$week9Due.Text    = "Assignment 8"
$week11Lesson.Text = "Data Wrangling IV"
$week11Due.Text    = "Reproducible Report: Initial Analysis"
